$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2331520.68
$ws.Range("C7").Value = -47.52465257126741
$ws.Range("D7").Value = 2367
$ws.Range("E7").Value = 2367
$ws.Range("F7").Value = 985.0108491761724
$ws.Range("G7").Value = 4.995033976543106
